# WithDraw Deal Testcase Update
#
# Adds 6 new automated test-case rows to the "Login" worksheet, right
# after the existing ShareDeal_ rows: DuplicateDeal_TC001/TC001(2)/TC002
# and WithDrawDeal_TC001/TC002/TC003. Each new row follows the same
# layout/formatting as the existing CreateDeal_/ShareDeal_ rows
# (Automation Test ID, UserName, Password, Expected Result), and the
# sheet's scroll/selection position is updated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

$stanEmail = "rogerdeals21+stan@gmail.com"
$johnEmail = "rogerdeals21+john@gmail.com"
$password  = "arewethere?"
$loginOk   = "Login successful"

$newRows = @(
    @{ Row = 26; Id = "DuplicateDeal_TC001";    User = $stanEmail },
    @{ Row = 27; Id = "DuplicateDeal_TC001(2)"; User = $johnEmail },
    @{ Row = 28; Id = "DuplicateDeal_TC002";    User = $stanEmail },
    @{ Row = 29; Id = "WithDrawDeal_TC001";     User = $stanEmail },
    @{ Row = 30; Id = "WithDrawDeal_TC002";     User = $stanEmail },
    @{ Row = 31; Id = "WithDrawDeal_TC003";     User = $johnEmail }
)

foreach ($r in $newRows) {
    # Copy the formatting of the last existing deal row (19) down onto the
    # new row so the new rows keep the same styling as CreateDeal_/ShareDeal_.
    $ws.Range("A19:D19").Copy($ws.Range("A" + $r.Row + ":D" + $r.Row))

    $ws.Cells.Item($r.Row, 1).Value2 = $r.Id
    $ws.Cells.Item($r.Row, 2).Value2 = $r.User
    $ws.Cells.Item($r.Row, 3).Value2 = $password
    $ws.Cells.Item($r.Row, 4).Value2 = $loginOk
}

# Update the sheet view to match the new scroll/selection position
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B21").Select()
